$wb = $excel.ActiveWorkbook

# --- Sheet: y_fitted_on_begin_2016 ---
$ws1 = $wb.Worksheets.Item("y_fitted_on_begin_2016")
$ws1.Cells.Item(2, 1).Value2 = 1981
$ws1.Cells.Item(2, 2).Value2 = 42.55384712233425
$ws1.Cells.Item(3, 1).Value2 = 1982
$ws1.Cells.Item(3, 2).Value2 = 43.78365725362809
$ws1.Cells.Item(4, 1).Value2 = 1983
$ws1.Cells.Item(4, 2).Value2 = 44.41638809225592
$ws1.Cells.Item(5, 1).Value2 = 1984
$ws1.Cells.Item(5, 2).Value2 = 45.34542409115377
$ws1.Cells.Item(6, 1).Value2 = 1985
$ws1.Cells.Item(6, 2).Value2 = 45.33563432462013
$ws1.Cells.Item(7, 1).Value2 = 1986
$ws1.Cells.Item(7, 2).Value2 = 43.91872048444355
$ws1.Cells.Item(8, 1).Value2 = 1987
$ws1.Cells.Item(8, 2).Value2 = 43.31497960742549
$ws1.Cells.Item(9, 1).Value2 = 1988
$ws1.Cells.Item(9, 2).Value2 = 43.90584469845007
$ws1.Cells.Item(10, 1).Value2 = 1989
$ws1.Cells.Item(10, 2).Value2 = 46.9042909337068
$ws1.Cells.Item(11, 1).Value2 = 1990
$ws1.Cells.Item(11, 2).Value2 = 49.93763364208019
$ws1.Cells.Item(12, 1).Value2 = 1991
$ws1.Cells.Item(12, 2).Value2 = 58.86859258404412
$ws1.Cells.Item(13, 1).Value2 = 1992
$ws1.Cells.Item(13, 2).Value2 = 60.17537149689453
$ws1.Cells.Item(14, 1).Value2 = 1993
$ws1.Cells.Item(14, 2).Value2 = 65.66365872099851
$ws1.Cells.Item(15, 1).Value2 = 1994
$ws1.Cells.Item(15, 2).Value2 = 70.16911639147581
$ws1.Cells.Item(16, 1).Value2 = 1995
$ws1.Cells.Item(16, 2).Value2 = 57.25802882685607
$ws1.Cells.Item(17, 1).Value2 = 1996
$ws1.Cells.Item(17, 2).Value2 = 58.0031298970992
$ws1.Cells.Item(18, 1).Value2 = 1997
$ws1.Cells.Item(18, 2).Value2 = 60.15090326350687
$ws1.Cells.Item(19, 1).Value2 = 1998
$ws1.Cells.Item(19, 2).Value2 = 60.4048114713287
$ws1.Cells.Item(20, 1).Value2 = 1999
$ws1.Cells.Item(20, 2).Value2 = 60.99550280958242
$ws1.Cells.Item(21, 1).Value2 = 2000
$ws1.Cells.Item(21, 2).Value2 = 60.28754235990747
$ws1.Cells.Item(22, 1).Value2 = 2001
$ws1.Cells.Item(22, 2).Value2 = 58.30732000718417
$ws1.Cells.Item(23, 1).Value2 = 2002
$ws1.Cells.Item(23, 2).Value2 = 58.44306496883486
$ws1.Cells.Item(24, 1).Value2 = 2003
$ws1.Cells.Item(24, 2).Value2 = 58.26887872988658
$ws1.Cells.Item(25, 1).Value2 = 2004
$ws1.Cells.Item(25, 2).Value2 = 57.43807271479668
$ws1.Cells.Item(26, 1).Value2 = 2005
$ws1.Cells.Item(26, 2).Value2 = 55.99438960970242
$ws1.Cells.Item(27, 1).Value2 = 2006
$ws1.Cells.Item(27, 2).Value2 = 56.94187965498509
$ws1.Cells.Item(28, 1).Value2 = 2007
$ws1.Cells.Item(28, 2).Value2 = 57.88472188355502
$ws1.Cells.Item(29, 1).Value2 = 2008
$ws1.Cells.Item(29, 2).Value2 = 58.45805016538195
$ws1.Cells.Item(30, 1).Value2 = 2009
$ws1.Cells.Item(30, 2).Value2 = 57.73080582122609
$ws1.Cells.Item(31, 1).Value2 = 2010
$ws1.Cells.Item(31, 2).Value2 = 59.66976718818956
$ws1.Cells.Item(32, 1).Value2 = 2011
$ws1.Cells.Item(32, 2).Value2 = 58.14571658755179
$ws1.Cells.Item(33, 1).Value2 = 2012
$ws1.Cells.Item(33, 2).Value2 = 58.30922071282652
$ws1.Cells.Item(34, 1).Value2 = 2013
$ws1.Cells.Item(34, 2).Value2 = 59.24002587400417
$ws1.Cells.Item(35, 1).Value2 = 2014
$ws1.Cells.Item(35, 2).Value2 = 59.93873815020661
$ws1.Cells.Item(36, 1).Value2 = 2015
$ws1.Cells.Item(36, 2).Value2 = 61.1785909343644
$ws1.Cells.Item(37, 1).Value2 = 2016
$ws1.Cells.Item(37, 2).Value2 = 62.00328041240942

# --- Sheet: y_pred_on_2017_2021 ---
$ws2 = $wb.Worksheets.Item("y_pred_on_2017_2021")
$ws2.Cells.Item(2, 2).Value2 = 62.70280038371784
$ws2.Cells.Item(3, 2).Value2 = 62.63005629749631
$ws2.Cells.Item(4, 2).Value2 = 62.66815399184938
$ws2.Cells.Item(5, 2).Value2 = 62.74589723294227
$ws2.Cells.Item(6, 2).Value2 = 62.84511013919229

# --- Sheet: y_fitted_on_begin_2021 ---
$ws3 = $wb.Worksheets.Item("y_fitted_on_begin_2021")
$ws3.Cells.Item(2, 1).Value2 = 1981
$ws3.Cells.Item(2, 2).Value2 = 42.5821526097135
$ws3.Cells.Item(3, 1).Value2 = 1982
$ws3.Cells.Item(3, 2).Value2 = 43.81213176070953
$ws3.Cells.Item(4, 1).Value2 = 1983
$ws3.Cells.Item(4, 2).Value2 = 44.44384922156458
$ws3.Cells.Item(5, 1).Value2 = 1984
$ws3.Cells.Item(5, 2).Value2 = 45.37175885521574
$ws3.Cells.Item(6, 1).Value2 = 1985
$ws3.Cells.Item(6, 2).Value2 = 45.36024403740812
$ws3.Cells.Item(7, 1).Value2 = 1986
$ws3.Cells.Item(7, 2).Value2 = 43.94059262082325
$ws3.Cells.Item(8, 1).Value2 = 1987
$ws3.Cells.Item(8, 2).Value2 = 43.33708486191995
$ws3.Cells.Item(9, 1).Value2 = 1988
$ws3.Cells.Item(9, 2).Value2 = 43.9273094544459
$ws3.Cells.Item(10, 1).Value2 = 1989
$ws3.Cells.Item(10, 2).Value2 = 46.92545453777105
$ws3.Cells.Item(11, 1).Value2 = 1990
$ws3.Cells.Item(11, 2).Value2 = 49.95787037913205
$ws3.Cells.Item(12, 1).Value2 = 1991
$ws3.Cells.Item(12, 2).Value2 = 58.85285158822337
$ws3.Cells.Item(13, 1).Value2 = 1992
$ws3.Cells.Item(13, 2).Value2 = 60.19607002537804
$ws3.Cells.Item(14, 1).Value2 = 1993
$ws3.Cells.Item(14, 2).Value2 = 65.65008469262791
$ws3.Cells.Item(15, 1).Value2 = 1994
$ws3.Cells.Item(15, 2).Value2 = 70.16520884744349
$ws3.Cells.Item(16, 1).Value2 = 1995
$ws3.Cells.Item(16, 2).Value2 = 57.2618169203752
$ws3.Cells.Item(17, 1).Value2 = 1996
$ws3.Cells.Item(17, 2).Value2 = 57.9907252205822
$ws3.Cells.Item(18, 1).Value2 = 1997
$ws3.Cells.Item(18, 2).Value2 = 60.14358121546393
$ws3.Cells.Item(19, 1).Value2 = 1998
$ws3.Cells.Item(19, 2).Value2 = 60.39024291901684
$ws3.Cells.Item(20, 1).Value2 = 1999
$ws3.Cells.Item(20, 2).Value2 = 60.97846616676454
$ws3.Cells.Item(21, 1).Value2 = 2000
$ws3.Cells.Item(21, 2).Value2 = 60.26347681024335
$ws3.Cells.Item(22, 1).Value2 = 2001
$ws3.Cells.Item(22, 2).Value2 = 58.26214194243693
$ws3.Cells.Item(23, 1).Value2 = 2002
$ws3.Cells.Item(23, 2).Value2 = 58.41107282395584
$ws3.Cells.Item(24, 1).Value2 = 2003
$ws3.Cells.Item(24, 2).Value2 = 58.23474116851791
$ws3.Cells.Item(25, 1).Value2 = 2004
$ws3.Cells.Item(25, 2).Value2 = 57.39674856467572
$ws3.Cells.Item(26, 1).Value2 = 2005
$ws3.Cells.Item(26, 2).Value2 = 55.93754633356588
$ws3.Cells.Item(27, 1).Value2 = 2006
$ws3.Cells.Item(27, 2).Value2 = 56.87927192000616
$ws3.Cells.Item(28, 1).Value2 = 2007
$ws3.Cells.Item(28, 2).Value2 = 57.81851360702088
$ws3.Cells.Item(29, 1).Value2 = 2008
$ws3.Cells.Item(29, 2).Value2 = 58.38908937185207
$ws3.Cells.Item(30, 1).Value2 = 2009
$ws3.Cells.Item(30, 2).Value2 = 57.656768311798
$ws3.Cells.Item(31, 1).Value2 = 2010
$ws3.Cells.Item(31, 2).Value2 = 59.59142584886293
$ws3.Cells.Item(32, 1).Value2 = 2011
$ws3.Cells.Item(32, 2).Value2 = 58.05240240194312
$ws3.Cells.Item(33, 1).Value2 = 2012
$ws3.Cells.Item(33, 2).Value2 = 58.21756772274924
$ws3.Cells.Item(34, 1).Value2 = 2013
$ws3.Cells.Item(34, 2).Value2 = 59.14323338498231
$ws3.Cells.Item(35, 1).Value2 = 2014
$ws3.Cells.Item(35, 2).Value2 = 59.83544747593455
$ws3.Cells.Item(36, 1).Value2 = 2015
$ws3.Cells.Item(36, 2).Value2 = 61.0704406387778
$ws3.Cells.Item(37, 1).Value2 = 2016
$ws3.Cells.Item(37, 2).Value2 = 61.89104513486561
$ws3.Cells.Item(38, 1).Value2 = 2017
$ws3.Cells.Item(38, 2).Value2 = 62.58824208419121
$ws3.Cells.Item(39, 1).Value2 = 2018
$ws3.Cells.Item(39, 2).Value2 = 62.71785594199115
$ws3.Cells.Item(40, 1).Value2 = 2019
$ws3.Cells.Item(40, 2).Value2 = 62.08910310231897
$ws3.Cells.Item(41, 1).Value2 = 2020
$ws3.Cells.Item(41, 2).Value2 = 62.4805715132816
$ws3.Cells.Item(42, 1).Value2 = 2021
$ws3.Cells.Item(42, 2).Value2 = 62.36777230784316
$ws3.Rows.Item(43).Delete() | Out-Null

# --- Sheet: y_pred_on_2022_2026 ---
$ws4 = $wb.Worksheets.Item("y_pred_on_2022_2026")
$ws4.Cells.Item(2, 2).Value2 = 60.5000686779338
$ws4.Cells.Item(3, 2).Value2 = 60.94167972307329
$ws4.Cells.Item(4, 2).Value2 = 61.18581539905068
$ws4.Cells.Item(5, 2).Value2 = 61.36690489210849
$ws4.Cells.Item(6, 2).Value2 = 61.52336745921818
